$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1011.98413
$ws.Range("I132").Value = 765.63464
$ws.Range("K132").Value = 2296.90392
$ws.Range("M132").Value = 233.0960800000003
$ws.Range("H137").Value = 1097.5676
$ws.Range("I137").Value = 901.2344000000001
$ws.Range("J137").Value = 2354.1
$ws.Range("K137").Value = 2703.7032
$ws.Range("L137").Value = 7062.299999999999
$ws.Range("M137").Value = -153.7031999999999
$ws.Range("N137").Value = -12162.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3327.2246
$ws.Range("I61").Value = 3839.8918
$ws.Range("J61").Value = 1746.5
$ws.Range("K61").Value = 3839.8918
$ws.Range("L61").Value = 1746.5
$ws.Range("M61").Value = -3627.8918
$ws.Range("N61").Value = -2170.5
$ws.Range("H110").Value = 1373.7273
$ws.Range("I110").Value = 838.875
$ws.Range("J110").Value = 2800
$ws.Range("K110").Value = 838.875
$ws.Range("L110").Value = 2800
$ws.Range("M110").Value = 1206.125
$ws.Range("N110").Value = -6890
$ws.Range("H132").Value = 1494575.1
$ws.Range("I132").Value = 1618.5454
$ws.Range("J132").Value = 4350666
$ws.Range("K132").Value = 4855.6362
$ws.Range("L132").Value = 13051998
$ws.Range("M132").Value = -2325.6362
$ws.Range("N132").Value = -13057058
$ws.Range("H136").Value = 3327.2246
$ws.Range("I136").Value = 3839.8918
$ws.Range("J136").Value = 1746.5
$ws.Range("K136").Value = 11519.6754
$ws.Range("L136").Value = 5239.5
$ws.Range("M136").Value = -8969.6754
$ws.Range("N136").Value = -10339.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 66667812
$ws.Range("I99").Value = 90909880
$ws.Range("J99").Value = 2125
$ws.Range("K99").Value = 90909880
$ws.Range("L99").Value = 2125
$ws.Range("M99").Value = -90908382
$ws.Range("N99").Value = -5121
$ws.Range("H107").Value = 83334350
$ws.Range("I107").Value = 125000930
$ws.Range("J107").Value = 1206.5
$ws.Range("K107").Value = 125000930
$ws.Range("L107").Value = 1206.5
$ws.Range("M107").Value = -124999010
$ws.Range("N107").Value = -5046.5
$ws.Range("H134").Value = 3190.7454
$ws.Range("I134").Value = 3379.8223
$ws.Range("K134").Value = 10139.4669
$ws.Range("M134").Value = -7604.466899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5391.4053
$ws.Range("I31").Value = 1327.7273
$ws.Range("J31").Value = 14704
$ws.Range("K31").Value = 1327.7273
$ws.Range("L31").Value = 14704
$ws.Range("M31").Value = -1032.7273
$ws.Range("N31").Value = -15294
$ws.Range("H34").Value = 5391.4053
$ws.Range("I34").Value = 1327.7273
$ws.Range("J34").Value = 14704
$ws.Range("K34").Value = 1327.7273
$ws.Range("L34").Value = 14704
$ws.Range("M34").Value = -1125.7273
$ws.Range("N34").Value = -15108
$ws.Range("H58").Value = 1313.2046
$ws.Range("I58").Value = 782.5484
$ws.Range("J58").Value = 2578.6155
$ws.Range("K58").Value = 782.5484
$ws.Range("L58").Value = 2578.6155
$ws.Range("M58").Value = -579.5484
$ws.Range("N58").Value = -2984.6155
$ws.Range("H125").Value = 17000
$ws.Range("J125").Value = 17000
$ws.Range("L125").Value = 17000
$ws.Range("N125").Value = -21920
$ws.Range("H132").Value = 1511.0702
$ws.Range("I132").Value = 1162.9
$ws.Range("J132").Value = 2330.2942
$ws.Range("K132").Value = 3488.7
$ws.Range("L132").Value = 6990.882599999999
$ws.Range("M132").Value = -958.7000000000003
$ws.Range("N132").Value = -12050.8826
$ws.Range("H136").Value = 1313.2046
$ws.Range("I136").Value = 782.5484
$ws.Range("J136").Value = 2578.6155
$ws.Range("K136").Value = 2347.6452
$ws.Range("L136").Value = 7735.8465
$ws.Range("M136").Value = 202.3548000000001
$ws.Range("N136").Value = -12835.8465

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 103808.9
$ws.Range("I5").Value = 186.27026
$ws.Range("J5").Value = 286382.1
$ws.Range("K5").Value = 558.81078
$ws.Range("L5").Value = 859146.2999999999
$ws.Range("M5").Value = -446.81078
$ws.Range("N5").Value = -859370.2999999999
$ws.Range("H70").Value = 2124.125
$ws.Range("I70").Value = 969.75
$ws.Range("J70").Value = 3278.5
$ws.Range("K70").Value = 2909.25
$ws.Range("L70").Value = 9835.5
$ws.Range("M70").Value = -2594.25
$ws.Range("N70").Value = -10465.5
$ws.Range("H73").Value = 2124.125
$ws.Range("I73").Value = 969.75
$ws.Range("J73").Value = 3278.5
$ws.Range("K73").Value = 2909.25
$ws.Range("L73").Value = 9835.5
$ws.Range("M73").Value = -1817.25
$ws.Range("N73").Value = -12019.5
$ws.Range("H135").Value = 103808.9
$ws.Range("I135").Value = 186.27026
$ws.Range("J135").Value = 286382.1
$ws.Range("K135").Value = 1676.43234
$ws.Range("L135").Value = 2577438.9
$ws.Range("M135").Value = 858.5676599999999
$ws.Range("N135").Value = -2582508.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 7490000
$ws.Range("I7").Value = 9250000
$ws.Range("J7").Value = 450000
$ws.Range("K7").Value = 9250000
$ws.Range("L7").Value = 450000
$ws.Range("M7").Value = -9249888
$ws.Range("N7").Value = -450224
$ws.Range("H8").Value = 7490000
$ws.Range("I8").Value = 9250000
$ws.Range("J8").Value = 450000
$ws.Range("K8").Value = 9250000
$ws.Range("L8").Value = 450000
$ws.Range("M8").Value = -9249861
$ws.Range("N8").Value = -450278
$ws.Range("H102").Value = 1357.375
$ws.Range("I102").Value = 1042.1
$ws.Range("J102").Value = 1882.8334
$ws.Range("K102").Value = 1042.1
$ws.Range("L102").Value = 1882.8334
$ws.Range("M102").Value = 579.9000000000001
$ws.Range("N102").Value = -5126.8334
$ws.Range("H126").Value = 5002.2
$ws.Range("I126").Value = 10614.728
$ws.Range("K126").Value = 31844.184
$ws.Range("M126").Value = -29374.184

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 24393156
$ws.Range("I40").Value = 32260226
$ws.Range("J40").Value = 5237
$ws.Range("K40").Value = 32260226
$ws.Range("L40").Value = 5237
$ws.Range("M40").Value = -32260090
$ws.Range("N40").Value = -5509
$ws.Range("H132").Value = 14066272
$ws.Range("I132").Value = 20557236
$ws.Range("J132").Value = 2516.3333
$ws.Range("K132").Value = 61671708
$ws.Range("L132").Value = 7548.999899999999
$ws.Range("M132").Value = -61669178
$ws.Range("N132").Value = -12608.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 18062.482
$ws.Range("I132").Value = 24448.357
$ws.Range("J132").Value = 1299.5625
$ws.Range("K132").Value = 73345.071
$ws.Range("L132").Value = 3898.6875
$ws.Range("M132").Value = -70815.071
$ws.Range("N132").Value = -8958.6875
$ws.Range("H136").Value = 6580957.5
$ws.Range("I136").Value = 2194.12
$ws.Range("J136").Value = 19232426
$ws.Range("K136").Value = 6582.36
$ws.Range("L136").Value = 57697278
$ws.Range("M136").Value = -4032.36
$ws.Range("N136").Value = -57702378
